# TC05_CDS_Filter_PHSAccession-phs001787.xlsx correction:
# The "FilesTab" query (cell B4 on Sheet1) filtered on
# experimental_strategies: ["RNA-Seq"] - remove that filter value so the
# query runs against all experimental strategies again (input file
# correction per commit "cds SCRIPTS 1-10 obj correction input file
# correction").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B4")
$query = $cell.Value()
$oldFilter = 'experimental_strategies: ["RNA-Seq"]'
$newFilter = 'experimental_strategies: []'
if ($query.Contains($oldFilter)) {
    $cell.Value = $query.Replace($oldFilter, $newFilter)

    # Setting the cell text retriggers Excel's wrap-text autofit on this
    # (wrapped) row; restore the row to its original rendered height.
    $ws.Rows.Item(4).RowHeight = 409.5
}

# Match the saved selection/viewport from the authored edit.
[void]$ws.Range("C4").Select()
